$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 203.5  # H2: 427.08334 -> 203.5
$ws.Cells.Item(2, 9).Value = 202.16667  # I2: 286.14285 -> 202.16667
$ws.Cells.Item(2, 10).Value = 205.5  # J2: 624.4 -> 205.5
$ws.Cells.Item(2, 11).Value = 202.16667  # K2: 286.14285 -> 202.16667
$ws.Cells.Item(2, 12).Value = 205.5  # L2: 624.4 -> 205.5
$ws.Cells.Item(2, 13).Value = -89.16667000000001  # M2: -173.14285 -> -89.16667000000001
$ws.Cells.Item(2, 14).Value = -431.5  # N2: -850.4 -> -431.5

# Row 28
$ws.Cells.Item(28, 8).Value = 1656.3334  # H28: 2445.9092 -> 1656.3334
$ws.Cells.Item(28, 9).Value = 613.375  # I28: 767.3333 -> 613.375
$ws.Cells.Item(28, 10).Value = 10000  # J28: 9999.5 -> 10000
$ws.Cells.Item(28, 11).Value = 613.375  # K28: 767.3333 -> 613.375
$ws.Cells.Item(28, 12).Value = 10000  # L28: 9999.5 -> 10000
$ws.Cells.Item(28, 13).Value = -128.375  # M28: -282.3333 -> -128.375
$ws.Cells.Item(28, 14).Value = -10970  # N28: -10969.5 -> -10970

# Row 98
$ws.Cells.Item(98, 8).Value = 1855.0278  # H98: 1806.081 -> 1855.0278
$ws.Cells.Item(98, 9).Value = 1856.6  # I98: 1806.25 -> 1856.6
$ws.Cells.Item(98, 11).Value = 1856.6  # K98: 1806.25 -> 1856.6
$ws.Cells.Item(98, 13).Value = -358.5999999999999  # M98: -308.25 -> -358.5999999999999

# Row 107
$ws.Cells.Item(107, 8).Value = 824  # H107: 858.93335 -> 824
$ws.Cells.Item(107, 9).Value = 836  # I107: 877.2308 -> 836
$ws.Cells.Item(107, 11).Value = 836  # K107: 877.2308 -> 836
$ws.Cells.Item(107, 13).Value = 1084  # M107: 1042.7692 -> 1084

# Row 122
$ws.Cells.Item(122, 8).Value = 1855.0278  # H122: 1806.081 -> 1855.0278
$ws.Cells.Item(122, 9).Value = 1856.6  # I122: 1806.25 -> 1856.6
$ws.Cells.Item(122, 11).Value = 5569.799999999999  # K122: 5418.75 -> 5569.799999999999
$ws.Cells.Item(122, 13).Value = -3119.799999999999  # M122: -2968.75 -> -3119.799999999999

# Row 135
$ws.Cells.Item(135, 8).Value = 710.3158  # H135: 729.1579 -> 710.3158
$ws.Cells.Item(135, 9).Value = 470.05554  # I135: 489.94446 -> 470.05554
$ws.Cells.Item(135, 11).Value = 4230.49986  # K135: 4409.50014 -> 4230.49986
$ws.Cells.Item(135, 13).Value = -1695.49986  # M135: -1874.50014 -> -1695.49986

# Row 137
$ws.Cells.Item(137, 8).Value = 3132371.5  # H137: 2506061.8 -> 3132371.5
$ws.Cells.Item(137, 9).Value = 6251297  # I137: 4167805.2 -> 6251297
$ws.Cells.Item(137, 11).Value = 18753891  # K137: 12503415.6 -> 18753891
$ws.Cells.Item(137, 13).Value = -18751341  # M137: -12500865.6 -> -18751341

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1519.8  # H2: 1449.8334 -> 1519.8
$ws.Cells.Item(2, 9).Value = 1333.3334  # I2: 1166.6666 -> 1333.3334
$ws.Cells.Item(2, 10).Value = 1799.5  # J2: 1733 -> 1799.5
$ws.Cells.Item(2, 11).Value = 1333.3334  # K2: 1166.6666 -> 1333.3334
$ws.Cells.Item(2, 12).Value = 1799.5  # L2: 1733 -> 1799.5
$ws.Cells.Item(2, 13).Value = -1220.3334  # M2: -1053.6666 -> -1220.3334
$ws.Cells.Item(2, 14).Value = -2025.5  # N2: -1959 -> -2025.5

# Row 32
$ws.Cells.Item(32, 8).Value = 3329.484  # H32: 3318.5247 -> 3329.484
$ws.Cells.Item(32, 10).Value = 4506.5  # J32: 4676 -> 4506.5
$ws.Cells.Item(32, 12).Value = 4506.5  # L32: 4676 -> 4506.5
$ws.Cells.Item(32, 14).Value = -5080.5  # N32: -5250 -> -5080.5

# Row 61
$ws.Cells.Item(61, 8).Value = 5297.952  # H61: 3617.8948 -> 5297.952
$ws.Cells.Item(61, 9).Value = 3682.9333  # I61: 2545.8438 -> 3682.9333
$ws.Cells.Item(61, 11).Value = 3682.9333  # K61: 2545.8438 -> 3682.9333
$ws.Cells.Item(61, 13).Value = -3470.9333  # M61: -2333.8438 -> -3470.9333

# Row 74
$ws.Cells.Item(74, 8).Value = 312374.06  # H74: 187571.67 -> 312374.06
$ws.Cells.Item(74, 9).Value = 2780285  # I74: 328022.75 -> 2780285
$ws.Cells.Item(74, 10).Value = 3885.1875  # J74: 3904.8462 -> 3885.1875
$ws.Cells.Item(74, 11).Value = 2780285  # K74: 328022.75 -> 2780285
$ws.Cells.Item(74, 12).Value = 3885.1875  # L74: 3904.8462 -> 3885.1875
$ws.Cells.Item(74, 13).Value = -2779411  # M74: -327148.75 -> -2779411
$ws.Cells.Item(74, 14).Value = -5633.1875  # N74: -5652.8462 -> -5633.1875

# Row 77
$ws.Cells.Item(77, 8).Value = 312374.06  # H77: 187571.67 -> 312374.06
$ws.Cells.Item(77, 9).Value = 2780285  # I77: 328022.75 -> 2780285
$ws.Cells.Item(77, 10).Value = 3885.1875  # J77: 3904.8462 -> 3885.1875
$ws.Cells.Item(77, 11).Value = 13901425  # K77: 1640113.75 -> 13901425
$ws.Cells.Item(77, 12).Value = 19425.9375  # L77: 19524.231 -> 19425.9375
$ws.Cells.Item(77, 13).Value = -13897057  # M77: -1635745.75 -> -13897057
$ws.Cells.Item(77, 14).Value = -28161.9375  # N77: -28260.231 -> -28161.9375

# Row 97
$ws.Cells.Item(97, 8).Value = 935.4  # H97: 996.6923 -> 935.4
$ws.Cells.Item(97, 9).Value = 990.2727  # I97: 996.13635 -> 990.2727
$ws.Cells.Item(97, 10).Value = 533  # J97: 999.75 -> 533
$ws.Cells.Item(97, 11).Value = 990.2727  # K97: 996.13635 -> 990.2727
$ws.Cells.Item(97, 12).Value = 533  # L97: 999.75 -> 533
$ws.Cells.Item(97, 13).Value = -494.2727  # M97: -500.13635 -> -494.2727
$ws.Cells.Item(97, 14).Value = -1525  # N97: -1991.75 -> -1525

# Row 110
$ws.Cells.Item(110, 8).Value = 3545.4783  # H110: 3114.2222 -> 3545.4783
$ws.Cells.Item(110, 9).Value = 1976.0714  # I110: 1677.9445 -> 1976.0714
$ws.Cells.Item(110, 11).Value = 1976.0714  # K110: 1677.9445 -> 1976.0714
$ws.Cells.Item(110, 13).Value = 68.92859999999996  # M110: 367.0554999999999 -> 68.92859999999996

# Row 116
$ws.Cells.Item(116, 8).Value = 1519.8  # H116: 1449.8334 -> 1519.8
$ws.Cells.Item(116, 9).Value = 1333.3334  # I116: 1166.6666 -> 1333.3334
$ws.Cells.Item(116, 10).Value = 1799.5  # J116: 1733 -> 1799.5
$ws.Cells.Item(116, 11).Value = 1333.3334  # K116: 1166.6666 -> 1333.3334
$ws.Cells.Item(116, 12).Value = 1799.5  # L116: 1733 -> 1799.5
$ws.Cells.Item(116, 13).Value = 960.6666  # M116: 1127.3334 -> 960.6666
$ws.Cells.Item(116, 14).Value = -6387.5  # N116: -6321 -> -6387.5

# Row 136
$ws.Cells.Item(136, 8).Value = 5297.952  # H136: 3617.8948 -> 5297.952
$ws.Cells.Item(136, 9).Value = 3682.9333  # I136: 2545.8438 -> 3682.9333
$ws.Cells.Item(136, 11).Value = 11048.7999  # K136: 7637.5314 -> 11048.7999
$ws.Cells.Item(136, 13).Value = -8498.7999  # M136: -5087.5314 -> -8498.7999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1519.8  # H3: 1449.8334 -> 1519.8
$ws.Cells.Item(3, 9).Value = 1333.3334  # I3: 1166.6666 -> 1333.3334
$ws.Cells.Item(3, 10).Value = 1799.5  # J3: 1733 -> 1799.5
$ws.Cells.Item(3, 11).Value = 1333.3334  # K3: 1166.6666 -> 1333.3334
$ws.Cells.Item(3, 12).Value = 1799.5  # L3: 1733 -> 1799.5
$ws.Cells.Item(3, 13).Value = -1219.3334  # M3: -1052.6666 -> -1219.3334
$ws.Cells.Item(3, 14).Value = -2027.5  # N3: -1961 -> -2027.5

# Row 94
$ws.Cells.Item(94, 8).Value = 142858050  # H94: 133334190 -> 142858050
$ws.Cells.Item(94, 9).Value = 153846370  # I94: 142857340 -> 153846370
$ws.Cells.Item(94, 11).Value = 153846370  # K94: 142857340 -> 153846370
$ws.Cells.Item(94, 13).Value = -153845919  # M94: -142856889 -> -153845919

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 1671.4286  # H22: 1916.6666 -> 1671.4286
$ws.Cells.Item(22, 9).Value = 1540.2  # I22: 1875.25 -> 1540.2
$ws.Cells.Item(22, 11).Value = 1540.2  # K22: 1875.25 -> 1540.2
$ws.Cells.Item(22, 13).Value = -1190.2  # M22: -1525.25 -> -1190.2

# Row 31
$ws.Cells.Item(31, 8).Value = 6029  # H31: 5441.6772 -> 6029
$ws.Cells.Item(31, 9).Value = 4319.636  # I31: 3561.6667 -> 4319.636
$ws.Cells.Item(31, 11).Value = 4319.636  # K31: 3561.6667 -> 4319.636
$ws.Cells.Item(31, 13).Value = -4024.636  # M31: -3266.6667 -> -4024.636

# Row 34
$ws.Cells.Item(34, 8).Value = 6029  # H34: 5441.6772 -> 6029
$ws.Cells.Item(34, 9).Value = 4319.636  # I34: 3561.6667 -> 4319.636
$ws.Cells.Item(34, 11).Value = 4319.636  # K34: 3561.6667 -> 4319.636
$ws.Cells.Item(34, 13).Value = -4117.636  # M34: -3359.6667 -> -4117.636

# Row 58
$ws.Cells.Item(58, 8).Value = 2361.7805  # H58: 2383.6829 -> 2361.7805
$ws.Cells.Item(58, 9).Value = 1780.4445  # I58: 1813.7037 -> 1780.4445
$ws.Cells.Item(58, 11).Value = 1780.4445  # K58: 1813.7037 -> 1780.4445
$ws.Cells.Item(58, 13).Value = -1577.4445  # M58: -1610.7037 -> -1577.4445

# Row 62
$ws.Cells.Item(62, 8).Value = 12511151  # H62: 10010070 -> 12511151
$ws.Cells.Item(62, 9).Value = 33338408  # I62: 25006304 -> 33338408
$ws.Cells.Item(62, 10).Value = 14797  # J62: 12580.833 -> 14797
$ws.Cells.Item(62, 11).Value = 33338408  # K62: 25006304 -> 33338408
$ws.Cells.Item(62, 12).Value = 14797  # L62: 12580.833 -> 14797
$ws.Cells.Item(62, 13).Value = -33337784  # M62: -25005680 -> -33337784
$ws.Cells.Item(62, 14).Value = -16045  # N62: -13828.833 -> -16045

# Row 65
$ws.Cells.Item(65, 8).Value = 12511151  # H65: 10010070 -> 12511151
$ws.Cells.Item(65, 9).Value = 33338408  # I65: 25006304 -> 33338408
$ws.Cells.Item(65, 10).Value = 14797  # J65: 12580.833 -> 14797
$ws.Cells.Item(65, 11).Value = 166692040  # K65: 125031520 -> 166692040
$ws.Cells.Item(65, 12).Value = 73985  # L65: 62904.165 -> 73985
$ws.Cells.Item(65, 13).Value = -166688920  # M65: -125028400 -> -166688920
$ws.Cells.Item(65, 14).Value = -80225  # N65: -69144.16500000001 -> -80225

# Row 110
$ws.Cells.Item(110, 8).Value = 0  # H110: 50704 -> 0
$ws.Cells.Item(110, 10).Value = 0  # J110: 50704 -> 0
$ws.Cells.Item(110, 12).Value = 0  # L110: 50704 -> 0
$ws.Cells.Item(110, 14).ClearContents() | Out-Null  # N110: -58884 -> (removed)

# Row 132
$ws.Cells.Item(132, 8).Value = 1826.8125  # H132: 1830.8959 -> 1826.8125
$ws.Cells.Item(132, 9).Value = 1574.1163  # I132: 1578.6744 -> 1574.1163
$ws.Cells.Item(132, 11).Value = 4722.3489  # K132: 4736.023200000001 -> 4722.3489
$ws.Cells.Item(132, 13).Value = -2192.3489  # M132: -2206.023200000001 -> -2192.3489

# Row 136
$ws.Cells.Item(136, 8).Value = 2361.7805  # H136: 2383.6829 -> 2361.7805
$ws.Cells.Item(136, 9).Value = 1780.4445  # I136: 1813.7037 -> 1780.4445
$ws.Cells.Item(136, 11).Value = 5341.333500000001  # K136: 5441.1111 -> 5341.333500000001
$ws.Cells.Item(136, 13).Value = -2791.333500000001  # M136: -2891.1111 -> -2791.333500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Cells.Item(8, 8).Value = 996.5  # H8: 997 -> 996.5
$ws.Cells.Item(8, 9).Value = 996.5  # I8: 997 -> 996.5
$ws.Cells.Item(8, 11).Value = 2989.5  # K8: 2991 -> 2989.5
$ws.Cells.Item(8, 13).Value = -2850.5  # M8: -2852 -> -2850.5

# Row 86
$ws.Cells.Item(86, 8).Value = 304.83334  # H86: 302.7143 -> 304.83334
$ws.Cells.Item(86, 9).Value = 310  # I86: 306.66666 -> 310
$ws.Cells.Item(86, 11).Value = 930  # K86: 919.9999799999999 -> 930
$ws.Cells.Item(86, 13).Value = 256  # M86: 266.0000200000001 -> 256

# Row 89
$ws.Cells.Item(89, 8).Value = 304.83334  # H89: 302.7143 -> 304.83334
$ws.Cells.Item(89, 9).Value = 310  # I89: 306.66666 -> 310
$ws.Cells.Item(89, 11).Value = 2790  # K89: 2759.99994 -> 2790
$ws.Cells.Item(89, 13).Value = 3138  # M89: 3168.00006 -> 3138

# Row 141
$ws.Cells.Item(141, 8).Value = 12330.692  # H141: 10783.866 -> 12330.692
$ws.Cells.Item(141, 9).Value = 4029.9  # I141: 3479.8333 -> 4029.9
$ws.Cells.Item(141, 11).Value = 12089.7  # K141: 10439.4999 -> 12089.7
$ws.Cells.Item(141, 13).Value = -6909.700000000001  # M141: -5259.499899999999 -> -6909.700000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 100021.91  # H70: 100026.86 -> 100021.91
$ws.Cells.Item(70, 9).Value = 129406.375  # I70: 129412.875 -> 129406.375
$ws.Cells.Item(70, 11).Value = 129406.375  # K70: 129412.875 -> 129406.375
$ws.Cells.Item(70, 13).Value = -129136.375  # M70: -129142.875 -> -129136.375

# Row 73
$ws.Cells.Item(73, 8).Value = 100021.91  # H73: 100026.86 -> 100021.91
$ws.Cells.Item(73, 9).Value = 129406.375  # I73: 129412.875 -> 129406.375
$ws.Cells.Item(73, 11).Value = 129406.375  # K73: 129412.875 -> 129406.375
$ws.Cells.Item(73, 13).Value = -128470.375  # M73: -128476.875 -> -128470.375

# Row 97
$ws.Cells.Item(97, 8).Value = 10000  # H97: 1383.619 -> 10000
$ws.Cells.Item(97, 9).Value = 0  # I97: 486.55554 -> 0
$ws.Cells.Item(97, 10).Value = 10000  # J97: 6766 -> 10000
$ws.Cells.Item(97, 11).Value = 0  # K97: 486.55554 -> 0
$ws.Cells.Item(97, 12).Value = 10000  # L97: 6766 -> 10000
$ws.Cells.Item(97, 13).ClearContents() | Out-Null  # M97: 9.444459999999992 -> (removed)
$ws.Cells.Item(97, 14).Value = -10992  # N97: -7758 -> -10992

# Row 113
$ws.Cells.Item(113, 8).Value = 2558.5557  # H113: 2746.8572 -> 2558.5557
$ws.Cells.Item(113, 9).Value = 2558.5557  # I113: 2746.8572 -> 2558.5557
$ws.Cells.Item(113, 11).Value = 2558.5557  # K113: 2746.8572 -> 2558.5557
$ws.Cells.Item(113, 13).Value = -388.5556999999999  # M113: -576.8571999999999 -> -388.5556999999999

# Row 139
$ws.Cells.Item(139, 8).Value = 81059.10000000001  # H139: 81359.10000000001 -> 81059.10000000001
$ws.Cells.Item(139, 10).Value = 82255  # J139: 82588.336 -> 82255
$ws.Cells.Item(139, 12).Value = 82255  # L139: 82588.336 -> 82255
$ws.Cells.Item(139, 14).Value = -92535  # N139: -92868.336 -> -92535

# Row 140
$ws.Cells.Item(140, 8).Value = 69998.55  # H140: 74460 -> 69998.55
$ws.Cells.Item(140, 10).Value = 69998.55  # J140: 74460 -> 69998.55
$ws.Cells.Item(140, 12).Value = 69998.55  # L140: 74460 -> 69998.55
$ws.Cells.Item(140, 14).Value = -80358.55  # N140: -84820 -> -80358.55

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 3021.6667  # H7: 3305.6 -> 3021.6667
$ws.Cells.Item(7, 9).Value = 3114.5454  # I7: 3450.6667 -> 3114.5454
$ws.Cells.Item(7, 11).Value = 3114.5454  # K7: 3450.6667 -> 3114.5454
$ws.Cells.Item(7, 13).Value = -3002.5454  # M7: -3338.6667 -> -3002.5454

# Row 126
$ws.Cells.Item(126, 8).Value = 3021.6667  # H126: 3305.6 -> 3021.6667
$ws.Cells.Item(126, 9).Value = 3114.5454  # I126: 3450.6667 -> 3114.5454
$ws.Cells.Item(126, 11).Value = 9343.636200000001  # K126: 10352.0001 -> 9343.636200000001
$ws.Cells.Item(126, 13).Value = -6873.636200000001  # M126: -7882.000100000001 -> -6873.636200000001

# Row 136
$ws.Cells.Item(136, 8).Value = 4621.294  # H136: 4804.3125 -> 4621.294
$ws.Cells.Item(136, 9).Value = 4547.9165  # I136: 4807.4546 -> 4547.9165
$ws.Cells.Item(136, 11).Value = 13643.7495  # K136: 14422.3638 -> 13643.7495
$ws.Cells.Item(136, 13).Value = -11093.7495  # M136: -11872.3638 -> -11093.7495
